$wb = $excel.ActiveWorkbook

# Sheet 1: FE_LFT_#1 - add row 56
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(56, 1).Value = 45842.49663194444
$ws.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 2).Value = "0x01,0x7c"
$ws.Cells.Item(56, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(56, 4).Value = "0x01,0x54"
$ws.Cells.Item(56, 5).Value = "0xf"
$ws.Cells.Item(56, 6).Value = 380
$ws.Cells.Item(56, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(56, 8).Value = 340
$ws.Cells.Item(56, 9).Value = 15

# Sheet 2: FE_LFT_#2 - add row 56
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(56, 1).Value = 45842.49663194444
$ws.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 2).Value = "0x01,0x90"
$ws.Cells.Item(56, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(56, 4).Value = "0x01,0x64"
$ws.Cells.Item(56, 5).Value = "0xe"
$ws.Cells.Item(56, 6).Value = 400
$ws.Cells.Item(56, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(56, 8).Value = 356
$ws.Cells.Item(56, 9).Value = 14

# Sheet 3: FE_PLT_#1 - add row 56
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(56, 1).Value = 45842.49663194444
$ws.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 2).Value = "0x00,0x6e"
$ws.Cells.Item(56, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(56, 4).Value = "0x00,0x68"
$ws.Cells.Item(56, 5).Value = "0x3"
$ws.Cells.Item(56, 6).Value = 110
$ws.Cells.Item(56, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(56, 8).Value = 104
$ws.Cells.Item(56, 9).Value = 3

# Sheet 4: FE_PLT_#2 - add row 56
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(56, 1).Value = 45842.49663194444
$ws.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 2).Value = "0x00,0x6e"
$ws.Cells.Item(56, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(56, 4).Value = "0x00,0x67"
$ws.Cells.Item(56, 5).Value = "0x3"
$ws.Cells.Item(56, 6).Value = 110
$ws.Cells.Item(56, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(56, 8).Value = 103
$ws.Cells.Item(56, 9).Value = 3
